# Update countries & provincias Spain
# Applies the data refresh captured by the diff:
#  - Updated "datos actualizados" timestamp in A1
#  - Updated case figures for a number of countries
#  - El Salvador's case count overtakes Kenia and Australia, so it moves
#    up one spot in the (descending, by total cases) ranking; Kenia and
#    Australia's rows shift down by one position (their own figures are
#    unchanged), while El Salvador gets the freshly reported numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 11:12"

# --- Helper to write a full data row (columns B..H) --------------------
function Set-Row($Row, $B, $C, $D, $E, $F, $G, $H) {
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
}

# --- Plain value refreshes (country/ranking unchanged) ------------------
Set-Row 4   3479573 90   1549624 1791702 0 0  138247   # Estados Unidos
Set-Row 6   910174  2529 573283  313112  0 52 23779    # India
Set-Row 20  190057  3163 103227  84406   0 33 2424     # Banglades
Set-Row 29  78572   1591 37636   37226   0 54 3710     # Indonesia
Set-Row 35  59568   1389 37987   21308   0 14 273      # Oman
Set-Row 36  57545   539  20459   35483   0 6  1603     # Filipinas
Set-Row 44  46630   347  42541   4063    0 0  26       # Singapur
Set-Row 46  41235   603  19474   21393   0 3  368      # Israel
Set-Row 47  38457   267  27756   9113    0 12 1588     # Polonia
Set-Row 118 1908    6    1493    387     0 0  28       # Eslovaquia
Set-Row 121 1875    1    1571    225     0 0  79       # Lituania
Set-Row 122 1859    10   1488    260     0 0  111      # Eslovenia

# --- El Salvador overtakes Kenia and Australia in the ranking ----------
# Row 73 (Sudan) and row 77 (Venezuela) keep their place and values.
$ws.Range("A74").Value = "El Salvador"
Set-Row 74 10303 325 5919 4106 0 11 278

$ws.Range("A75").Value = "Kenia"
Set-Row 75 10294 0   2946 7151 0 0  197

$ws.Range("A76").Value = "Australia"
Set-Row 76 10250 270 7835 2307 0 0  108
